# Penalty Reward System edit (unfinished, per commit message).
#
# "Weekly Quantity" sheet: remove the week-of-2024-03-10 row (old row 24,
# A24=45361.99999999999 / B24=847) and shift all following rows up by one.
# "Monthly Trend" sheet: row 9 (A9=45382.99999999999) requested quantity
# changes from 2548 to 1701.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(24).Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B9").Value = 1701
